$d = $word.ActiveDocument

# 1. Azioni Utente (row 1): replace the customer action description.
#    "1." stays untouched; the rest becomes the new sentence about viewing the cart.
$d.Content.Find.Execute(
    "Il Cliente inserisce prodotti nel carrello indicandone la quantità.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Il Cliente decide di voler visualizzare il proprio carrello.", 2) | Out-Null

# 2. Azioni Sistema (row 1): renumber from 3. to 2. and change the system action text.
$d.Content.Find.Execute(
    "3.Il Sistema registra le azioni del Cliente all’interno del carrello.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2.Il Sistema fa visualizzare il carrello al Cliente.", 2) | Out-Null

# 3. Exit Condition, paragraph 1: "L'utente ha smesso..." -> "Il Cliente prosegue con l'ordinazione."
$d.Content.Find.Execute(
    "L’utente ha smesso di apportare modifiche ai prodotti nel carrello.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Il Cliente prosegue con l'ordinazione.", 2) | Out-Null

# 4. Exit Condition, paragraph 2: "L'utente prosegue con l'ordinazione." -> "Il Cliente non è più loggato nel sistema."
$d.Content.Find.Execute(
    "L’utente prosegue con l’ordinazione.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Il Cliente non è più loggato nel sistema.", 2) | Out-Null

# 5. Exit Condition, paragraph 3: "L'utente non è più loggato nel sistema." -> "Il Cliente esce dalla visualizzazione."
$d.Content.Find.Execute(
    "L’utente non è più loggato nel sistema.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Il Cliente esce dalla visualizzazione.", 2) | Out-Null

# The "_GoBack" bookmark (last-edit marker) now belongs at the end of the
# Exit Condition text (end of the paragraph we just edited) rather than after
# the old Azioni Utente sentence.
$lastPara = $d.Paragraphs(29).Range
$bmPos = $lastPara.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null
